$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, columns B, C, D, E, G (F unchanged)
$data = @{
    2 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 }
    3 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 8.660232485948974;   G = 13.71653804550039 }
    4 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 8.660232485948974;   G = 14.36450238910742 }
    5 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    6 = @{ B = 0.6753301551942219; C = 0.3127903958511391; D = 0.1575252929769615; E = 8.660232485948974;   G = 9.805878329971296 }
    7 = @{ B = 0.127881588408715;  C = 0.3127903958511391; D = 3.900430680208489;  E = 0.496779210170732;  G = 4.837881874639075 }
    8 = @{ B = 0.04763786555579896; C = 0.3127903958511391; D = 0.1575252929769615; E = 0.496779210170732; G = 1.014732764554632 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
